$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.060.23'
$ws.Range("E2").Value = '  +3.68%  '
$ws.Range("D3").Value = '1.726.80'
$ws.Range("D5").Value = '218.97'
$ws.Range("E5").Value = '  +1.66%  '
$ws.Range("D6").Value = '0.523'
$ws.Range("E6").Value = '  +1.25%  '
$ws.Range("D8").Value = '24.11'
$ws.Range("E8").Value = '  +13.28%  '
$ws.Range("E9").Value = '  +3.33%  '
$ws.Range("D10").Value = '0.0632'
$ws.Range("E10").Value = '  +1.86%  '
$ws.Range("D11").Value = '0.0901'
$ws.Range("E11").Value = '  +1.90%  '
$ws.Range("D13").Value = '1.730.80'
$ws.Range("E13").Value = '  +3.21%  '
$ws.Range("E14").Value = '  +3.38%  '
$ws.Range("E15").Value = '  +5.41%  '
$ws.Range("D16").Value = '67.61'
$ws.Range("E16").Value = '  +2.50%  '
$ws.Range("D17").Value = '28.025.66'
$ws.Range("E17").Value = '  +3.57%  '
$ws.Range("D18").Value = '243.19'
$ws.Range("E18").Value = '  +2.43%  '
$ws.Range("D19").Value = '0.0₃0754'
$ws.Range("E19").Value = '  +1.84%  '
$ws.Range("E20").Value = '  -3.19%  '
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("E22").Value = '  +3.82%  '
$ws.Range("D23").Value = '9.73'
$ws.Range("E23").Value = '  +4.30%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("D25").Value = '148.96'
$ws.Range("E25").Value = '  +1.48%  '
$ws.Range("D26").Value = '7.52'
$ws.Range("E26").Value = '  +4.36%  '
$ws.Range("D27").Value = '16.74'
$ws.Range("E27").Value = '  +2.42%  '
$ws.Range("D28").Value = '0.114'
$ws.Range("E28").Value = '  +1.47%  '
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("E30").Value = '  +2.49%  '
$ws.Range("E31").Value = '  +2.04%  '
$ws.Range("E32").Value = '  +2.82%  '
$ws.Range("D33").Value = '1.495.58'
$ws.Range("E33").Value = '  -3.23%  '
$ws.Range("D34").Value = '3.27'
$ws.Range("E34").Value = '  +2.55%  '
$ws.Range("E35").Value = '  -2.44%  '
$ws.Range("D36").Value = '0.954'
$ws.Range("E36").Value = '  +3.29%  '
$ws.Range("D37").Value = '0.606'
$ws.Range("E37").Value = '  +1.11%  '
$ws.Range("E38").Value = '  +0.67%  '
$ws.Range("E39").Value = '  +0.41%  '
$ws.Range("D40").Value = '1.07'
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("D41").Value = '70.84'
$ws.Range("E41").Value = '  +4.79%  '
$ws.Range("D42").Value = '5.83'
$ws.Range("E42").Value = '  +4.18%  '
$ws.Range("E44").Value = '  +2.30%  '
$ws.Range("D45").Value = '1.874.69'
$ws.Range("E45").Value = '  +2.83%  '
$ws.Range("D46").Value = '0.796'
$ws.Range("E46").Value = '  +1.89%  '
$ws.Range("E47").Value = '  +12.03%  '
$ws.Range("D48").Value = '91.18'
$ws.Range("E48").Value = '  +0.55%  '
$ws.Range("E49").Value = '  +3.63%  '
$ws.Range("E50").Value = '  +0.81%  '
$ws.Range("E51").Value = '  +1.14%  '
